# Update TPM-derived values for the "ECs" sending/target cluster.
# The underlying ligand (Lgals3bp) and receptor (Itgb1) average/total
# expression values for the ECs cluster were recomputed with new TPM
# data. That change cascades into the specificity columns (I/J, O/P)
# and the edge-weight columns (Q/R, S/T) for every row that references
# ECs as either the Sending cluster or the Target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New base TPM values (ECs cluster) ---
$G_ECs = 16.52389666666667
$H_ECs = 49.57169
$M_ECs = 77.08952333333333
$N_ECs = 231.26857

# --- Unchanged base TPM values for the other clusters ---
$G_FAPs  = 30.06295833333333
$H_FAPs  = 90.188875
$M_FAPs  = 101.5800373333333
$N_FAPs  = 304.740112

$G_MuSCs = 5.611374666666666
$H_MuSCs = 16.834124
$M_MuSCs = 142.0267893333333
$N_MuSCs = 426.080368

# Sums across the three clusters, used for specificity (I/J, O/P)
$sumG = $G_ECs + $G_FAPs + $G_MuSCs
$sumH = $H_ECs + $H_FAPs + $H_MuSCs
$sumM = $M_ECs + $M_FAPs + $M_MuSCs
$sumN = $N_ECs + $N_FAPs + $N_MuSCs

# Row layout: (row, sendingCluster, targetCluster, whether G/H (sending) or
# M/N (target) changed for this row -> determines whether Q/R need
# recomputation or can be left at their original values)
$rows = @(
  @{ Row = 2;  G = $G_ECs;   H = $H_ECs;   M = $M_ECs;   N = $N_ECs;   Recalc = $true  },
  @{ Row = 3;  G = $G_ECs;   H = $H_ECs;   M = $M_FAPs;  N = $N_FAPs;  Recalc = $true  },
  @{ Row = 4;  G = $G_ECs;   H = $H_ECs;   M = $M_MuSCs; N = $N_MuSCs; Recalc = $true  },
  @{ Row = 5;  G = $G_FAPs;  H = $H_FAPs;  M = $M_ECs;   N = $N_ECs;   Recalc = $true  },
  @{ Row = 6;  G = $G_FAPs;  H = $H_FAPs;  M = $M_FAPs;  N = $N_FAPs;  Recalc = $false },
  @{ Row = 7;  G = $G_FAPs;  H = $H_FAPs;  M = $M_MuSCs; N = $N_MuSCs; Recalc = $false },
  @{ Row = 8;  G = $G_MuSCs; H = $H_MuSCs; M = $M_ECs;   N = $N_ECs;   Recalc = $true  },
  @{ Row = 9;  G = $G_MuSCs; H = $H_MuSCs; M = $M_FAPs;  N = $N_FAPs;  Recalc = $false },
  @{ Row = 10; G = $G_MuSCs; H = $H_MuSCs; M = $M_MuSCs; N = $N_MuSCs; Recalc = $false }
)

# First pass: write G, H, M, N, and compute I, J, O, P (specificities).
# Q (=G*M) and R (=H*N) are only rewritten for rows whose underlying G or M
# (resp. H or N) actually changed; otherwise the original Q/R values
# already on the sheet are reused (kept) to avoid float re-round noise.
$Qs = @{}
$Rs = @{}
foreach ($item in $rows) {
    $r = $item.Row
    $G = $item.G
    $H = $item.H
    $M = $item.M
    $N = $item.N

    $I = $G / $sumG
    $J = $H / $sumH
    $O = $M / $sumM
    $P = $N / $sumN

    $ws.Cells.Item($r, 7).Value  = $G   # G
    $ws.Cells.Item($r, 8).Value  = $H   # H
    $ws.Cells.Item($r, 9).Value  = $I   # I
    $ws.Cells.Item($r, 10).Value = $J   # J
    $ws.Cells.Item($r, 13).Value = $M   # M
    $ws.Cells.Item($r, 14).Value = $N   # N
    $ws.Cells.Item($r, 15).Value = $O   # O
    $ws.Cells.Item($r, 16).Value = $P   # P

    if ($item.Recalc) {
        $Q = $G * $M
        $R = $H * $N
        $ws.Cells.Item($r, 17).Value = $Q   # Q
        $ws.Cells.Item($r, 18).Value = $R   # R
        $Qs[$r] = $Q
        $Rs[$r] = $R
    }
    else {
        # Keep the existing (unchanged) Q/R values
        $Qs[$r] = $ws.Cells.Item($r, 17).Value2
        $Rs[$r] = $ws.Cells.Item($r, 18).Value2
    }
}

# Second pass: S and T are Q and R normalized against the sum of all rows.
# Because several Q/R values changed, every row's S/T must be refreshed
# even when that row's own Q/R value stayed the same.
$sumQ = 0
$sumR = 0
foreach ($item in $rows) {
    $r = $item.Row
    $sumQ += $Qs[$r]
    $sumR += $Rs[$r]
}

foreach ($item in $rows) {
    $r = $item.Row
    $S = $Qs[$r] / $sumQ
    $T = $Rs[$r] / $sumR
    $ws.Cells.Item($r, 19).Value = $S   # S
    $ws.Cells.Item($r, 20).Value = $T   # T
}
